# Update "想去人数" (want-to-go count) figures in F column on both the
# "展览" and "全部类型" worksheets to reflect refreshed counts from the
# regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - F column updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 189
$wsExpo.Range("F3").Value  = 232
$wsExpo.Range("F5").Value  = 790
$wsExpo.Range("F6").Value  = 240
$wsExpo.Range("F7").Value  = 5830
$wsExpo.Range("F8").Value  = 30
$wsExpo.Range("F11").Value = 44
$wsExpo.Range("F14").Value = 179
$wsExpo.Range("F15").Value = 340
$wsExpo.Range("F16").Value = 26

# Sheet "全部类型" (all types) - F column updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 189
$wsAll.Range("F3").Value  = 232
$wsAll.Range("F5").Value  = 790
$wsAll.Range("F6").Value  = 241
$wsAll.Range("F7").Value  = 5830
$wsAll.Range("F8").Value  = 30
$wsAll.Range("F11").Value = 44
$wsAll.Range("F14").Value = 179
$wsAll.Range("F15").Value = 340
$wsAll.Range("F16").Value = 26

$wb.Save()
